$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that gain the "new data" yellow highlight (style 48 -> 49) ---
# Copy formatting from a stable donor cell that already carries style 49
# (N13 keeps its own style 49 throughout this edit) onto C46:C51, then set
# their new "Latest Period" values.
$ws.Range("N13").Copy()
$ws.Range("C46:C51").PasteSpecial(-4122)

$ws.Range("C46").Value = 45992
$ws.Range("C47").Value = 45992
$ws.Range("C48").Value = 45992
$ws.Range("C49").Value = 45992
$ws.Range("C50").Value = 45992
$ws.Range("C51").Value = 45992

# --- Cell that loses the highlight (style 49 -> 48) ---
# Copy formatting from a stable donor cell that already carries style 48
# (C3 is untouched by this edit) onto N51; its value is unchanged (46062).
$ws.Range("C3").Copy()
$ws.Range("N51").PasteSpecial(-4122)
$ws.Range("N51").Value = 46062

# --- Row 13: UI Initial Claims ---
$ws.Range("N13").Value = 46062
$ws.Range("Q13").Value = 206000
$ws.Range("R13").Value = 229000
$ws.Range("S13").Value = 232000
$ws.Range("T13").Value = 209000
$ws.Range("U13").Value = 210000

# --- Row 14: UI Continuing Claims ---
$ws.Range("N14").Value = 46055
$ws.Range("Q14").Value = 1869000
$ws.Range("R14").Value = 1852000
$ws.Range("S14").Value = 1841000
$ws.Range("T14").Value = 1819000
$ws.Range("U14").Value = 1865000

# --- Row 29: 5yr, 5yr Forward ---
$ws.Range("N29").Value = 46071
$ws.Range("Q29").Value = 2.15
$ws.Range("R29").Value = 2.13
$ws.Range("S29").Value = 2.12
$ws.Range("T29").Value = 2.13
$ws.Range("U29").Value = 2.15

# --- Row 30: 10yr TIPS ---
$ws.Range("N30").Value = 46071
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.26
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.32

# --- Row 46: Exports (X), BOPTEXP level ---
$ws.Range("F46").Value = 287287
$ws.Range("G46").Value = 292290
$ws.Range("H46").Value = 302594
$ws.Range("I46").Value = 293901
$ws.Range("J46").Value = 283736

# --- Row 47: BOPTEXP M/M % Delta SA, and FFR ---
$ws.Range("F47").Value = -0.01711656231824554
$ws.Range("G47").Value = -0.03405222839844813
$ws.Range("H47").Value = 0.02957798714533122
$ws.Range("I47").Value = 0.03582555615078808
$ws.Range("J47").Value = 0.0004866042778863822
$ws.Range("N47").Value = 46070

# --- Row 48: Imports (M), BOPTIMP level, and 2y UST ---
$ws.Range("F48").Value = 357598
$ws.Range("G48").Value = 345334
$ws.Range("H48").Value = 331343
$ws.Range("I48").Value = 341582
$ws.Range("J48").Value = 338909
$ws.Range("N48").Value = 46070
$ws.Range("Q48").Value = 3.43
$ws.Range("R48").Value = 3.4
$ws.Range("S48").Value = 3.47
$ws.Range("T48").Value = 3.52
$ws.Range("U48").Value = 3.45

# --- Row 49: BOPTIMP M/M % Delta SA, and 5y UST ---
$ws.Range("F49").Value = 0.0355134449547394
$ws.Range("G49").Value = 0.04222512622871166
$ws.Range("H49").Value = -0.02997523288697879
$ws.Range("I49").Value = 0.007887072931081818
$ws.Range("J49").Value = -0.05210885495329198
$ws.Range("N49").Value = 46070
$ws.Range("Q49").Value = 3.63
$ws.Range("R49").Value = 3.61
$ws.Range("S49").Value = 3.67
$ws.Range("T49").Value = 3.75
$ws.Range("U49").Value = 3.7

# --- Row 50: Trade Balance (BOPSTB) level, and 10y UST ---
$ws.Range("F50").Value = 29018
$ws.Range("G50").Value = 30597
$ws.Range("H50").Value = 29777
$ws.Range("I50").Value = 30169
$ws.Range("J50").Value = 30416
$ws.Range("N50").Value = 46070
$ws.Range("Q50").Value = 4.05
$ws.Range("R50").Value = 4.04
$ws.Range("S50").Value = 4.09
$ws.Range("T50").Value = 4.18
$ws.Range("U50").Value = 4.16

# --- Row 51: BOPSTB M/M % Delta SA, and 30y Mortgage (highlight removed) ---
$ws.Range("F51").Value = -0.05160636663725204
$ws.Range("G51").Value = 0.02753803270980959
$ws.Range("H51").Value = -0.01299347011833341
$ws.Range("I51").Value = -0.008120725933719042
$ws.Range("J51").Value = 0.06327343913864225

# --- Row 52: BAA ---
$ws.Range("N52").Value = 46070
$ws.Range("Q52").Value = 5.75
$ws.Range("R52").Value = 5.76
$ws.Range("S52").Value = 5.77
$ws.Range("T52").Value = 5.85
$ws.Range("U52").Value = 5.82
